$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.098.16"
$ws.Range("E2").Value = "  +1.14%  "

$ws.Range("D3").Value = "'2.305.42"
$ws.Range("E3").Value = "  +0.72%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'300.75"
$ws.Range("E5").Value = "  -0.19%  "

$ws.Range("D6").Value = "'97.48"
$ws.Range("E6").Value = "  -0.87%  "

$ws.Range("D7").Value = "'0.509"
$ws.Range("E7").Value = "  +0.41%  "

$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "'0.503"
$ws.Range("E9").Value = "  -0.12%  "

$ws.Range("D10").Value = "'33.83"
$ws.Range("E10").Value = "  -2.05%  "

$ws.Range("D11").Value = "'0.0795"
$ws.Range("E11").Value = "  +0.27%  "

$ws.Range("D12").Value = "'49.39"
$ws.Range("E12").Value = "  -2.82%  "

$ws.Range("E13").Value = "  +2.52%  "

$ws.Range("D14").Value = "'17.19"
$ws.Range("E14").Value = "  +11.30%  "

$ws.Range("D15").Value = "'6.79"
$ws.Range("E15").Value = "  +0.86%  "

$ws.Range("D16").Value = "'2.663.54"
$ws.Range("E16").Value = "  +0.85%  "

$ws.Range("D17").Value = "'2.294.84"
$ws.Range("E17").Value = "  +0.63%  "

$ws.Range("D18").Value = "'0.815"
$ws.Range("E18").Value = "  +2.77%  "

$ws.Range("D19").Value = "'43.027.86"
$ws.Range("E19").Value = "  +1.18%  "

$ws.Range("D20").Value = "'11.68"
$ws.Range("E20").Value = "  +0.48%  "

$ws.Range("D21").Value = "'0.0₃0902"
$ws.Range("E21").Value = "  +0.52%  "

$ws.Range("D22").Value = "'6.09"
$ws.Range("E22").Value = "  +0.89%  "

$ws.Range("D23").Value = "'68.03"
$ws.Range("E23").Value = "  +1.39%  "

$ws.Range("D24").Value = "'237.12"
$ws.Range("E24").Value = "  +0.68%  "

$ws.Range("D25").Value = "'2.04"
$ws.Range("E25").Value = "  +4.26%  "

$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("D27").Value = "'2.46"
$ws.Range("E27").Value = "  -1.85%  "

$ws.Range("D28").Value = "'24.51"
$ws.Range("E28").Value = "  -0.59%  "

$ws.Range("E29").Value = "  -5.50%  "

$ws.Range("D30").Value = "'167.35"
$ws.Range("E30").Value = "  +1.49%  "

$ws.Range("D31").Value = "'34.21"
$ws.Range("E31").Value = "  -0.24%  "

$ws.Range("D32").Value = "'9.15"
$ws.Range("E32").Value = "  +0.16%  "

$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  -0.04%  "

$ws.Range("D34").Value = "'4.69"
$ws.Range("E34").Value = "  +7.35%  "

$ws.Range("D35").Value = "'4.98"
$ws.Range("E35").Value = "  -0.28%  "

$ws.Range("D37").Value = "'16.91"
$ws.Range("E37").Value = "  +4.29%  "

$ws.Range("D38").Value = "'0.0701"
$ws.Range("E38").Value = "  +0.06%  "

$ws.Range("E39").Value = "  -0.12%  "

$ws.Range("E40").Value = "  +0.71%  "

$ws.Range("D41").Value = "'1.78"
$ws.Range("E41").Value = "  -0.57%  "

$ws.Range("E42").Value = "  -0.23%  "

$ws.Range("E43").Value = "  +0.43%  "

$ws.Range("D44").Value = "'1.981.98"
$ws.Range("E44").Value = "  +0.59%  "

$ws.Range("D45").Value = "'0.0283"
$ws.Range("E45").Value = "  -0.17%  "

$ws.Range("D46").Value = "'9.90"
$ws.Range("E46").Value = "  +1.43%  "

$ws.Range("D47").Value = "'17.64"
$ws.Range("E47").Value = "  -2.78%  "

$ws.Range("D48").Value = "'2.87"
$ws.Range("E48").Value = "  +0.08%  "

$ws.Range("D49").Value = "'53.49"
$ws.Range("E49").Value = "  +0.58%  "

$ws.Range("D50").Value = "'2.528.56"
$ws.Range("E50").Value = "  +0.57%  "

$ws.Range("D51").Value = "'4.60"
